# Additional file tweaks from Obvience.
#
# The authoring/content-management tooling ("Obvience") re-synced the
# SharePoint content-type metadata that PowerPoint keeps in the package's
# custom XML parts:
#   - customXml/item1.xml (+ itemProps1.xml) -> ct:contentTypeSchema,
#     whose ma:versionID / ma:fieldsID were refreshed, and whose
#     datastore GUID (itemProps1.xml / CustomXMLPart.Id) was reissued.
#   - customXml/item2.xml (+ itemProps2.xml) and
#     customXml/item3.xml (+ itemProps3.xml) -> their datastore GUIDs
#     (itemProps2.xml / itemProps3.xml, i.e. each part's .Id) were
#     likewise reissued.
#
# Each CustomXMLPart in PowerPoint's object model corresponds to one
# customXml/itemN.xml payload (.XML) plus its datastore GUID, stored in
# customXml/itemPropsN.xml and surfaced as the part's .Id.

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

function Update-CustomXmlPart([int]$partIndex, [string]$oldId, [string]$newId, [string[]]$xmlReplacements) {
    if ($partIndex -lt 1 -or $partIndex -gt $parts.Count) { return }
    $part = $parts.Item($partIndex)

    if ($xmlReplacements.Count -gt 0) {
        $xml = $part.XML
        for ($i = 0; $i -lt $xmlReplacements.Count; $i += 2) {
            $xml = $xml.Replace($xmlReplacements[$i], $xmlReplacements[$i + 1])
        }
        $part.XML = $xml
    }

    if ($part.Id -eq $oldId) {
        $part.Id = $newId
    }
}

# customXml/item1.xml + itemProps1.xml
Update-CustomXmlPart 1 `
    "{6AA53FB1-E571-4D8C-A271-F26160660FEB}" `
    "{C1EDE165-DB7A-402C-AE07-CFB039D2DDA9}" `
    @(
        "4257d9291635bc5cb4c092963aa1be7a", "ade2e80b87cb37ba6f94e81c3365b849",
        "bca45b24ba46885b01554d6d0eb772ad", "ff2a7fbf0c0a09b394f21031231ef037"
    )

# customXml/item2.xml + itemProps2.xml
Update-CustomXmlPart 2 `
    "{A67FED21-8961-4231-8A56-C3BF8FC35199}" `
    "{ECF4DBA5-C131-4C8B-B191-B1C4F5F38522}" `
    @()

# customXml/item3.xml + itemProps3.xml
Update-CustomXmlPart 3 `
    "{215F13FB-9260-43F0-A978-9B04E29EAFB4}" `
    "{8F4BAE5D-2E17-4C02-9672-C0FC52F33382}" `
    @()
